# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.086.78'
$ws.Range('E2').Value = '  -3.82%  '
$ws.Range('D3').Value = '2.235.94'
$ws.Range('E3').Value = '  -4.77%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''232.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.06%  '
$ws.Range('D6').Value = '''0.636'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.42%  '
$ws.Range('D7').Value = '''71.35'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.20%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.568'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.34%  '
$ws.Range('D10').Value = '''0.0996'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('D11').Value = '''58.52'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').Value = '''35.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.26%  '
$ws.Range('E13').Value = '  -3.25%  '
$ws.Range('D14').Value = '''6.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.53%  '
$ws.Range('D15').Value = '2.569.82'
$ws.Range('E15').Value = '  -4.72%  '
$ws.Range('D16').Value = '''15.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.19%  '
$ws.Range('D17').Value = '''0.870'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.70%  '
$ws.Range('D18').Value = '2.248.22'
$ws.Range('E18').Value = '  -4.10%  '
$ws.Range('D19').Value = '42.000.23'
$ws.Range('E19').Value = '  -3.88%  '
$ws.Range('D20').Value = '0.0₃0988'
$ws.Range('E20').Value = '  -2.90%  '
$ws.Range('D21').Value = '''73.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.52%  '
$ws.Range('D22').Value = '''6.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.24%  '
$ws.Range('D23').Value = '''236.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.55%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '''1.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('D26').Value = '''3.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.28%  '
$ws.Range('E27').Value = '  -5.14%  '
$ws.Range('D28').Value = '''9.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.78%  '
$ws.Range('D29').Value = '''2.11'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.25%  '
$ws.Range('D30').Value = '''166.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.01%  '
$ws.Range('D31').Value = '''20.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.82%  '
$ws.Range('E32').Value = '  -7.24%  '
$ws.Range('E33').Value = '  -7.46%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '''5.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.0715'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.32%  '
$ws.Range('D36').Value = '''4.82'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.79%  '
$ws.Range('E37').Value = '  -6.08%  '
$ws.Range('D38').Value = '''22.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +17.51%  '
$ws.Range('D39').Value = '''6.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('D40').Value = '''2.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.38%  '
$ws.Range('D41').Value = '''0.0265'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.90%  '
$ws.Range('D42').Value = '''66.55'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('D43').Value = '''5.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('D44').Value = '''8.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.41%  '
$ws.Range('D45').Value = '''0.100'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.69%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('B47').Value = 'SynthetixNetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D47').Value = '''4.56'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +11.66%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '''0.187'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.41%  '
$ws.Range('D49').Value = '''2.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.00%  '
$ws.Range('E50').Value = '  -5.42%  '
$ws.Range('D51').Value = '''2.81'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.74%  '
